$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows above the current row 105 (pushes old rows 105-121 down to 111-127)
$ws.Range("A105:T110").EntireRow.Insert()

$newRows = @(
    @(2,"Comercializadora del Agro de Limarí","Coquimbo",44559,4,"Fruta",100103,"Frutos de hueso (carozo)",100103001,"Cereza","Lapins","Especial",500,9500,10000,9750,"`$/bandeja 10 kilos","Región de O'Higgins",975,10),
    @(2,"Comercializadora del Agro de Limarí","Coquimbo",44559,4,"Fruta",100103,"Frutos de hueso (carozo)",100103001,"Cereza","Lapins","Primera",400,7500,8000,7750,"`$/bandeja 10 kilos","Región de O'Higgins",775,10),
    @(2,"Comercializadora del Agro de Limarí","Coquimbo",44559,4,"Fruta",100103,"Frutos de hueso (carozo)",100103001,"Cereza","Lapins","Segunda",300,5500,6000,5750,"`$/bandeja 10 kilos","Región de O'Higgins",575,10),
    @(2,"Comercializadora del Agro de Limarí","Coquimbo",44559,4,"Fruta",100103,"Frutos de hueso (carozo)",100103001,"Cereza","Santina","Especial",500,9500,10000,9750,"`$/bandeja 10 kilos","Región de O'Higgins",975,10),
    @(2,"Comercializadora del Agro de Limarí","Coquimbo",44559,4,"Fruta",100103,"Frutos de hueso (carozo)",100103001,"Cereza","Santina","Primera",400,7500,8000,7750,"`$/bandeja 10 kilos","Región de O'Higgins",775,10),
    @(2,"Comercializadora del Agro de Limarí","Coquimbo",44559,4,"Fruta",100103,"Frutos de hueso (carozo)",100103001,"Cereza","Santina","Segunda",280,5500,6000,5750,"`$/bandeja 10 kilos","Región de O'Higgins",575,10)
)

for ($r = 0; $r -lt $newRows.Length; $r++) {
    $row = $newRows[$r]
    $targetRow = 105 + $r
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value = $row[$c]
    }
}
